{"js": "const replacements = [\n  [\"930\u00d73=2790\", \"964\u00d74=3856\"],\n  [\"350\u00d74=1400\", \"629\u00d75=3145\"],\n  [\"990\u00d72=1980\", \"208\u00d74=832\"],\n  [\"926\u00d76=5556\", \"297\u00d72=594\"],\n  [\"605\u00d74=2420\", \"159\u00d77=1113\"],\n  [\"280\u00d79=2520\", \"801\u00d79=7209\"],\n  [\"587\u00d74=2348\", \"435\u00d72=870\"],\n  [\"110\u00d75=550\", \"461\u00d72=922\"],\n  [\"674\u00d75=3370\", \"993\u00d73=2979\"],\n  [\"837\u00d73=2511\", \"850\u00d72=1700\"],\n  [\"499\u00d72=998\", \"895\u00d74=3580\"],\n  [\"584\u00d77=4088\", \"758\u00d79=6822\"],\n  [\"631\u00d77=4417\", \"305\u00d74=1220\"],\n  [\"889\u00d76=5334\", \"851\u00d72=1702\"],\n  [\"476\u00d77=3332\", \"167\u00d73=501\"],\n  [\"558\u00d74=2232\", \"355\u00d78=2840\"],\n  [\"363\u00d78=2904\", \"371\u00d75=1855\"],\n  [\"774\u00d76=4644\", \"963\u00d77=6741\"],\n  [\"135\u00d76=810\", \"134\u00d78=1072\"],\n  [\"600\u00d74=2400\", \"244\u00d73=732\"],\n  [\"977\u00d73=2931\", \"879\u00d72=1758\"],\n  [\"525\u00d76=3150\", \"410\u00d79=3690\"],\n  [\"517\u00d78=4136\", \"677\u00d74=2708\"],\n  [\"702\u00d72=1404\", \"378\u00d76=2268\"],\n  [\"188\u00d77=1316\", \"830\u00d74=3320\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  ,@(\"930\u00d73=2790\", \"964\u00d74=3856\")\n  ,@(\"350\u00d74=1400\", \"629\u00d75=3145\")\n  ,@(\"990\u00d72=1980\", \"208\u00d74=832\")\n  ,@(\"926\u00d76=5556\", \"297\u00d72=594\")\n  ,@(\"605\u00d74=2420\", \"159\u00d77=1113\")\n  ,@(\"280\u00d79=2520\", \"801\u00d79=7209\")\n  ,@(\"587\u00d74=2348\", \"435\u00d72=870\")\n  ,@(\"110\u00d75=550\", \"461\u00d72=922\")\n  ,@(\"674\u00d75=3370\", \"993\u00d73=2979\")\n  ,@(\"837\u00d73=2511\", \"850\u00d72=1700\")\n  ,@(\"499\u00d72=998\", \"895\u00d74=3580\")\n  ,@(\"584\u00d77=4088\", \"758\u00d79=6822\")\n  ,@(\"631\u00d77=4417\", \"305\u00d74=1220\")\n  ,@(\"889\u00d76=5334\", \"851\u00d72=1702\")\n  ,@(\"476\u00d77=3332\", \"167\u00d73=501\")\n  ,@(\"558\u00d74=2232\", \"355\u00d78=2840\")\n  ,@(\"363\u00d78=2904\", \"371\u00d75=1855\")\n  ,@(\"774\u00d76=4644\", \"963\u00d77=6741\")\n  ,@(\"135\u00d76=810\", \"134\u00d78=1072\")\n  ,@(\"600\u00d74=2400\", \"244\u00d73=732\")\n  ,@(\"977\u00d73=2931\", \"879\u00d72=1758\")\n  ,@(\"525\u00d76=3150\", \"410\u00d79=3690\")\n  ,@(\"517\u00d78=4136\", \"677\u00d74=2708\")\n  ,@(\"702\u00d72=1404\", \"378\u00d76=2268\")\n  ,@(\"188\u00d77=1316\", \"830\u00d74=3320\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n}"}
